$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.352.42'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +6.27%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.336.19'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.74%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '411.91'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '110.88'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.80%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.585'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.59%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.634'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.42'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0998'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.16%  '
$ws.Range('E12').Value = '  +1.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.859.26'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.31%  '
$ws.Range('E14').Value = '  +2.82%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '19.67'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.347.43'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.72%  '
$ws.Range('E17').Value = '  +1.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '60.095.91'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.98%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.74'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.55%  '
$ws.Range('E20').Value = '  +2.48%  '
$ws.Range('E21').Value = '  +5.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.12'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.84%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '301.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.79'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.18'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.28%  '
$ws.Range('B26').Value = 'RenderToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.90'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +9.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '28.62'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.59%  '
$ws.Range('B28').Value = 'LEO'
$ws.Range('C28').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '4.48'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.16%  '
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.03'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.28%  '
$ws.Range('E30').Value = '  +5.62%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.66'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +24.70%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.115'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.49'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.45%  '
$ws.Range('E34').Value = '  +0.37%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '39.46'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.82%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0506'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.35'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.76%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.11'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.32%  '
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('E40').Value = '  -3.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '137.77'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.96%  '
$ws.Range('E42').Value = '  +3.20%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.295'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.37%  '
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('E45').Value = '  -0.79%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.89'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.18%  '
$ws.Range('E47').Value = '  +8.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.53'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.85%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.179.94'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.94%  '
$ws.Range('E50').Value = '  +1.56%  '
$ws.Range('E51').Value = '  -0.61%  '
